$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: username@domain moves to a new address, password column matches the
# standard "password" label (previously held the literal "pass1234").
$ws.Range("A2").Value = "hicamod348@halbov.com"
$ws.Range("B2").Value = "password"

# Row 3 keeps its existing email / password text (shared-string reindex only).
$ws.Range("A3").Value = "payeba8662@baxima.com"
$ws.Range("B3").Value = "password"

# New hyperlink on A2 pointing at the new email address.
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:hicamod348@halbov.com") | Out-Null
$ws.Range("A2").Style = "Normal"
$wb.Styles.Item("Hyperlink").Delete()

# Selection moves to C2.
$ws.Range("C2").Select() | Out-Null
